$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (TC 001) updates ---
# B3: hyperlink cell text changes from old email to the new one; keep the
# existing hyperlink relationship/style untouched.
$ws.Range("B3").Value = "Haleluyaa.aki@gamail.com"

# C3: was a plain number (123456); becomes a text value that starts with
# "@" so it must be stored as text.
$ws.Range("C3").Value = "@Chai4704183;"

# --- Row 4 (TC 002) updates ---
# B4: was "u4704183"; becomes the phone number "0831539901" which must stay
# text (no leading-zero loss). Keep it formatted as text first so Excel
# doesn't coerce it to a number.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0831539901"

# C4: was the literal quoted text "0831539901"; becomes plain text "u4704183".
$ws.Range("C4").Value = "u4704183"

# D4 (Thanawit Chaisuphapsirikun) is unchanged.

# --- New Row 5 (TC 003) ---
# Copy A4's number format/alignment (right/center, text format) down to A5
# before writing its value so "003" keeps its leading zeros.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = "003"

$ws.Range("B5").Value = "nuthey@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:nuthey@hotmail.com") | Out-Null

$ws.Range("C5").Value = "Vasan247086;"
$ws.Range("D5").Value = "Natty Nattha Tualek"

# Match the saved selection/active cell.
$ws.Range("D13").Select() | Out-Null
